# Scheduled market-data refresh: updates currentAveragePrice(NQ/HQ)
# and the derived LevePrice/LeveProfit columns (H:N) for the rows
# whose item prices moved since the last snapshot, across all job sheets.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H113").Value = 19401.2
$ws.Range("I113").Value = 0
$ws.Range("J113").Value = 19401.2
$ws.Range("K113").Value = 0
$ws.Range("L113").Value = 19401.2
$ws.Range("M113").ClearContents()
$ws.Range("N113").Value = -25909.2
$ws.Range("H116").Value = 636312.8
$ws.Range("I116").Value = 1252624.1
$ws.Range("J116").Value = 20001.5
$ws.Range("K116").Value = 1252624.1
$ws.Range("L116").Value = 20001.5
$ws.Range("M116").Value = -1249182.1
$ws.Range("N116").Value = -26885.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H8").Value = 8816.5
$ws.Range("J8").Value = 10378.8
$ws.Range("L8").Value = 10378.8
$ws.Range("N8").Value = -10666.8
$ws.Range("H19").Value = 27500
$ws.Range("I19").Value = 25000
$ws.Range("J19").Value = 35000
$ws.Range("K19").Value = 25000
$ws.Range("L19").Value = 35000
$ws.Range("M19").Value = -24771
$ws.Range("N19").Value = -35458
$ws.Range("H32").Value = 11713.043
$ws.Range("I32").Value = 8098.9067
$ws.Range("K32").Value = 8098.9067
$ws.Range("M32").Value = -7811.9067
$ws.Range("H45").Value = 1858.6666
$ws.Range("I45").Value = 1022.4
$ws.Range("K45").Value = 1022.4
$ws.Range("M45").Value = -645.4
$ws.Range("H132").Value = 2099.1428
$ws.Range("I132").Value = 1018.5
$ws.Range("J132").Value = 5557.2
$ws.Range("K132").Value = 3055.5
$ws.Range("L132").Value = 16671.6
$ws.Range("M132").Value = -525.5
$ws.Range("N132").Value = -21731.6
$ws.Range("H137").Value = 53106.332
$ws.Range("J137").Value = 53106.332
$ws.Range("L137").Value = 53106.332
$ws.Range("N137").Value = -63306.332

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H59").Value = 52999.5
$ws.Range("J59").Value = 52999.5
$ws.Range("L59").Value = 52999.5
$ws.Range("N59").Value = -54693.5
$ws.Range("H134").Value = 2539.7273
$ws.Range("I134").Value = 1327.091
$ws.Range("J134").Value = 7390.273
$ws.Range("K134").Value = 3981.273
$ws.Range("L134").Value = 22170.819
$ws.Range("M134").Value = -1446.273
$ws.Range("N134").Value = -27240.819
$ws.Range("H135").Value = 54427.25
$ws.Range("J135").Value = 54427.25
$ws.Range("L135").Value = 54427.25
$ws.Range("N135").Value = -64567.25
$ws.Range("H137").Value = 35309.875
$ws.Range("J137").Value = 35309.875
$ws.Range("L137").Value = 35309.875
$ws.Range("N137").Value = -45509.875

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2795.9778
$ws.Range("I31").Value = 1313.1875
$ws.Range("J31").Value = 6445.923
$ws.Range("K31").Value = 1313.1875
$ws.Range("L31").Value = 6445.923
$ws.Range("M31").Value = -1018.1875
$ws.Range("N31").Value = -7035.923
$ws.Range("H34").Value = 2795.9778
$ws.Range("I34").Value = 1313.1875
$ws.Range("J34").Value = 6445.923
$ws.Range("K34").Value = 1313.1875
$ws.Range("L34").Value = 6445.923
$ws.Range("M34").Value = -1111.1875
$ws.Range("N34").Value = -6849.923
$ws.Range("H62").Value = 45460090
$ws.Range("J62").Value = 5768.875
$ws.Range("L62").Value = 5768.875
$ws.Range("N62").Value = -7016.875
$ws.Range("H65").Value = 45460090
$ws.Range("J65").Value = 5768.875
$ws.Range("L65").Value = 28844.375
$ws.Range("N65").Value = -35084.375
$ws.Range("H134").Value = 5704.654
$ws.Range("I134").Value = 6407.9473
$ws.Range("J134").Value = 3795.7144
$ws.Range("K134").Value = 19223.8419
$ws.Range("L134").Value = 11387.1432
$ws.Range("M134").Value = -16688.8419
$ws.Range("N134").Value = -16457.1432
$ws.Range("H139").Value = 48980
$ws.Range("J139").Value = 48980
$ws.Range("L139").Value = 48980
$ws.Range("N139").Value = -59260

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H26").Value = 7572.3
$ws.Range("I26").Value = 15193.5
$ws.Range("J26").Value = 2491.5
$ws.Range("K26").Value = 45580.5
$ws.Range("L26").Value = 7474.5
$ws.Range("M26").Value = -45292.5
$ws.Range("N26").Value = -8050.5
$ws.Range("H64").Value = 5513.8335
$ws.Range("J64").Value = 10653
$ws.Range("L64").Value = 31959
$ws.Range("N64").Value = -32499
$ws.Range("H67").Value = 5513.8335
$ws.Range("J67").Value = 10653
$ws.Range("L67").Value = 31959
$ws.Range("N67").Value = -33831
$ws.Range("H107").Value = 51078.4
$ws.Range("J107").Value = 127162.875
$ws.Range("L107").Value = 381488.625
$ws.Range("N107").Value = -385328.625
$ws.Range("H114").Value = 3709.3333
$ws.Range("H123").Value = 3999.6667
$ws.Range("J123").Value = 4999
$ws.Range("L123").Value = 14997
$ws.Range("N123").Value = -19897
$ws.Range("H137").Value = 3203.077
$ws.Range("I137").Value = 2919.0908
$ws.Range("J137").Value = 4765
$ws.Range("K137").Value = 8757.2724
$ws.Range("L137").Value = 14295
$ws.Range("M137").Value = -3657.2724
$ws.Range("N137").Value = -24495

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H137").Value = 77443.8
$ws.Range("J137").Value = 77443.8
$ws.Range("L137").Value = 77443.8
$ws.Range("N137").Value = -87643.8

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H10").Value = 90000
$ws.Range("J10").Value = 90000
$ws.Range("L10").Value = 90000
$ws.Range("N10").Value = -90280
$ws.Range("H40").Value = 7551.091
$ws.Range("I40").Value = 10112.6
$ws.Range("J40").Value = 5416.5
$ws.Range("K40").Value = 10112.6
$ws.Range("L40").Value = 5416.5
$ws.Range("M40").Value = -9976.6
$ws.Range("N40").Value = -5688.5
$ws.Range("H132").Value = 4870.1816
$ws.Range("I132").Value = 1893.1818
$ws.Range("J132").Value = 7847.1816
$ws.Range("K132").Value = 5679.5454
$ws.Range("L132").Value = 23541.5448
$ws.Range("M132").Value = -3149.5454
$ws.Range("N132").Value = -28601.5448

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H9").Value = 0
$ws.Range("I9").Value = 0
$ws.Range("J9").Value = 0
$ws.Range("K9").Value = 0
$ws.Range("L9").Value = 0
$ws.Range("M9").ClearContents()
$ws.Range("N9").ClearContents()
$ws.Range("H80").Value = 35757.285
$ws.Range("J80").Value = 35757.285
$ws.Range("L80").Value = 35757.285
$ws.Range("N80").Value = -37753.285
$ws.Range("H83").Value = 35757.285
$ws.Range("J83").Value = 35757.285
$ws.Range("L83").Value = 107271.855
$ws.Range("N83").Value = -117255.855
$ws.Range("H132").Value = 6670279
$ws.Range("I132").Value = 3934.2666
$ws.Range("J132").Value = 16669797
$ws.Range("K132").Value = 11802.7998
$ws.Range("L132").Value = 50009391
$ws.Range("M132").Value = -9272.799800000001
$ws.Range("N132").Value = -50014451
$ws.Range("H136").Value = 4828.3687
$ws.Range("I136").Value = 1583.9
$ws.Range("J136").Value = 8433.333000000001
$ws.Range("K136").Value = 4751.700000000001
$ws.Range("L136").Value = 25299.999
$ws.Range("M136").Value = -2201.700000000001
$ws.Range("N136").Value = -30399.999
